$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (A2:F7) covering the full operator comparison table.
$data = @(
    @("A1", "Axis",      2.5,  3.5,  2.5,  3),
    @("A2", "Indosat",   3.33, 2.83, 3.17, 3.33),
    @("A3", "Smartfren", 3,    3.5,  2,    3.5),
    @("A4", "Telkomsel", 4,    3.9,  4.7,  4),
    @("A5", "Tri",       3.67, 3.83, 2.5,  4),
    @("A6", "XL",        2.25, 2.5,  3.5,  2.75)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
